$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 99

# Write the date as literal text (matches source data where dates are
# stored as plain strings, not date serials) by temporarily forcing a
# text number format, then resetting the cell style back to the default
# so no extra formatting is left behind on the new cell.
$cell = $ws.Cells.Item($row, 1)
$cell.NumberFormat = "@"
$cell.Value = "11/24/2025"
$cell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 8258.360000000001
